$d = $word.ActiveDocument

# --- Change 1: Image1 (Shapes collection index 1 maps correctly) ---
# wp:wrapSquare -> wp:wrapTopAndBottom  (WrapFormat.Type 0 -> 4 in this host)
$d.Shapes.Item(1).WrapFormat.Type = 4

# --- Change 2: a new empty paragraph is inserted right before the paragraph
# that holds the "Image17" picture (currently Paragraphs.Item(7)). ---
$d.Paragraphs.Item(7).Range.InsertParagraphBefore()

# --- Change 3: Image17 repositioning + wrap mode.
# NOTE: this host's Shapes collection has an indexing quirk when *setting*
# properties (confirmed empirically) -- Shapes.Item(3) is the call that
# actually reaches the last shape ("Image17") for property assignment,
# while reads via Shapes.Item(17) correctly report "Image17" by name. ---
$img17 = $d.Shapes.Item(3)
$img17.Left = 3.75
$img17.Top = 0.8
$img17.WrapFormat.Type = 4
